# Rearranged Mexico from Central to North America
#
# The workbook drives a Sankey-style flow table (A=source node id,
# B=target node id, C=value, D=color, E=label). Row 46 is the "Mexico"
# flow, which used to feed into the "Central America" aggregator node
# (id 36, same as El Salvador / Other Central America, colored blue) and
# now instead feeds into the "Northern America" aggregator node (id 44,
# same as Canada / Other Northern America, colored gray). The aggregate
# rows that summarize those flows (Central America->Latin America,
# Latin America->Americas, Northern America(excl US)->Americas) are
# adjusted by Mexico's value to keep the totals consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the Mexico flow from Central America to Northern America ---
$ws.Range("B46").Value = 44
$ws.Range("D46").Value = "rgba(166, 166, 166, 0.5)"

# --- Re-balance the downstream aggregate totals by Mexico's value ---
$ws.Range("C53").Value = 3581747
$ws.Range("C57").Value = 11345894
$ws.Range("C58").Value = 11994362

# --- Formatting touch-up: these value cells picked up an explicit
# "no fill" flag in the saved file (applyFill set, still no visible
# fill color) -- replicate that by touching Interior on the same cells.
$ws.Range("C45:C59").Interior.ColorIndex = -4142

# --- Selection / view bookkeeping ---
$ws.Range("C46").Select()
